$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "W.A.T.N.Jayathilake"
$ws.Range("B2").Value = "Not Completed"

$ws.Range("A3").Value = "D.W.S.N.Sewwandi"
$ws.Range("B3").Value = "Not Completed"

$ws.Range("A4").Value = "L.R.M.U.BANDARA"
$ws.Range("B4").Value = "Not Completed"
